$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $text) {
    # Leading apostrophe forces Excel to store the value as literal text
    # (matches the workbook's existing convention of keeping numeric-looking
    # quantities/prices as text). Resetting the style afterwards strips the
    # "quote prefix" number-format flag that the apostrophe entry adds, so
    # the cell ends up as plain text with the sheet's default style - same
    # as every other cell in this sheet.
    $ws.Range($addr).Value = "'" + $text
    $ws.Range($addr).Style = "Normal"
}

# Row 22
$ws.Range("A22").Value = "DV"
$ws.Range("B22").Value = "Chia Seeds - Black"
Set-TextCell "C22" "1"
Set-TextCell "D22" "125.65"
Set-TextCell "E22" "125.65"

# Row 23
$ws.Range("A23").Value = "Palmer/Sysco"
$ws.Range("B23").Value = "Jam - Strawberry"
Set-TextCell "C23" "1"
Set-TextCell "D23" "0.00"
Set-TextCell "E23" "0.00"

# Row 24 (SKU column blank, same pattern as several earlier rows)
Set-TextCell "A24" ""
$ws.Range("B24").Value = "Mustard - Honey"
Set-TextCell "C24" "1"
Set-TextCell "D24" "0.00"
Set-TextCell "E24" "0.00"

# Row 25
Set-TextCell "A25" ""
$ws.Range("B25").Value = "Nuts - Walnut Halves & Pieces"
Set-TextCell "C25" "1"
Set-TextCell "D25" "3.08"
Set-TextCell "E25" "3.08"

# Row 26
Set-TextCell "A26" ""
$ws.Range("B26").Value = "Flour - Millers Choice"
Set-TextCell "C26" "1"
Set-TextCell "D26" "0.00"
Set-TextCell "E26" "0.00"
